$d = $word.ActiveDocument

function Merge-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

$rsq = [char]0x2019   # right single quotation mark (')
$endash = [char]0x2013   # en dash (-)

$t1 = "drivers of malnutrition."
Merge-Text $t1 $t1

$t2 = "major heal" + "th issue"
$t2r = "major health issue"
Merge-Text $t2 $t2r

$t3 = "Logi" + "stic Regression"
$t3r = "Logistic Regression"
Merge-Text $t3 $t3r

$t4 = "more affected than urban" + " peers"
$t4r = "more affected than urban peers"
Merge-Text $t4 $t4r

$t5 = "Socioeconomic" + ": Poverty"
$t5r = "Socioeconomic: Poverty"
Merge-Text $t5 $t5r

$t6 = "Water trucki" + "ng, handwashing"
$t6r = "Water trucking, handwashing"
Merge-Text $t6 $t6r

$t7 = "(3" + $endash + "5+ Years)"
Merge-Text $t7 $t7

$t8 = "central to Rwa" + "nda" + $rsq + "s success"
Merge-Text $t8 $t8

$t9 = "Engage" + " local leaders and women" + $rsq + "s groups"
Merge-Text $t9 $t9

$t10 = "multisectoral ac" + "tion, a healthier"
$t10r = "multisectoral action, a healthier"
Merge-Text $t10 $t10r

$t11 = "(2021). R" + "wanda Demographic"
$t11r = "(2021). Rwanda Demographic"
Merge-Text $t11 $t11r

$t12 = "Hidden Hunger Rwanda D" + "ashboard"
Merge-Text $t12 $t12

Write-Output "done"
